$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "reviews_count" column (E) is empty for every row and is being
# removed entirely, shifting columns F:K left to E:J.
$ws.Range("E:E").Delete()
